# Update the handoff/handback datetime stamps on the zh-cn and de-de
# report sheets, as produced by a re-run of the handback status report
# generation ("Generate Report for Handback").
#
# Correspond Handoff Datetime  -> column E (rows 2 and 4)
# Correspond Handback DateTime -> column H (rows 2 and 4)
#
# These cells already carry the "yyyy-mm-dd HH:mm:ss"-looking text, and
# are stored as plain shared-string text (not real date serials) in the
# original workbook, so we just overwrite the text in place.

$wb = $excel.ActiveWorkbook

$sheetUpdates = @(
    @{ Sheet = "zh-cn"; Cell = "E2"; Value = "2016-03-19 20:18:12" },
    @{ Sheet = "zh-cn"; Cell = "H2"; Value = "2016-03-19 20:18:31" },
    @{ Sheet = "zh-cn"; Cell = "E4"; Value = "2016-03-19 20:18:12" },
    @{ Sheet = "zh-cn"; Cell = "H4"; Value = "2016-03-19 20:18:31" },
    @{ Sheet = "de-de"; Cell = "E2"; Value = "2016-03-19 20:18:15" },
    @{ Sheet = "de-de"; Cell = "H2"; Value = "2016-03-19 20:18:36" },
    @{ Sheet = "de-de"; Cell = "E4"; Value = "2016-03-19 20:18:15" },
    @{ Sheet = "de-de"; Cell = "H4"; Value = "2016-03-19 20:18:36" }
)

foreach ($update in $sheetUpdates) {
    $ws = $wb.Worksheets.Item($update.Sheet)
    $range = $ws.Range($update.Cell)
    $range.Value = $update.Value
}
